$d = $word.ActiveDocument

# --- Paula Castellanos (Rol Líder paragraph): spelling fixes ---
$d.Content.Find.Execute("estuviramos", $true, $false, $false, $false, $false, $true, 1, $false, "estuviéramos", 2)
$d.Content.Find.Execute("usaramos", $true, $false, $false, $false, $false, $true, 1, $false, "usáramos", 2)
$d.Content.Find.Execute("cumpliento", $true, $false, $false, $false, $false, $true, 1, $false, "cumplimento", 2)

# --- Ingrid Echavarria -> Echavarría ---
$d.Content.Find.Execute("Ingrid Echavarria", $true, $false, $false, $false, $false, $true, 1, $false, "Ingrid Echavarría", 2)

# --- David Abril role (Rol Calidad paragraph): Tubo -> Tuvo, estandares -> estándares ---
$d.Content.Find.Execute("Tubo muy buen", $true, $false, $false, $false, $false, $true, 1, $false, "Tuvo muy buen", 2)
$d.Content.Find.Execute("documentos de estandares", $true, $false, $false, $false, $false, $true, 1, $false, "documentos de estándares", 2)

# --- Gabriel Martinez (Rol Ingeniero paragraph): "la momento" -> "al momento" ---
$d.Content.Find.Execute("muy detallista la momento", $true, $false, $false, $false, $false, $true, 1, $false, "muy detallista al momento", 2)

# --- Felipe Fagua (Rol Desarrollo paragraph): lider -> líder, demas -> demás ---
$d.Content.Find.Execute("como lider de desarrollo", $true, $false, $false, $false, $false, $true, 1, $false, "como líder de desarrollo", 2)
$d.Content.Find.Execute("opiniones de los demas integrantes", $true, $false, $false, $false, $false, $true, 1, $false, "opiniones de los demás integrantes", 2)

# --- Paula Castellanos (Autoevaluación, Rol Procesos paragraph): lider -> líder ---
$d.Content.Find.Execute("apoyo al lider desarrollo", $true, $false, $false, $false, $false, $true, 1, $false, "apoyo al líder desarrollo", 2)

# --- metodologias -> metodologías, mas detallista -> más detallista ---
$d.Content.Find.Execute("las metodologias empleadas", $true, $false, $false, $false, $false, $true, 1, $false, "las metodologías empleadas", 2)
$d.Content.Find.Execute("hizo falta ser mas detallista", $true, $false, $false, $false, $false, $true, 1, $false, "hizo falta ser más detallista", 2)

# --- Paula Castellanos (Autoevaluación, Rol Ingeniero paragraph): mas cuidadosa -> más cuidadosa ---
$d.Content.Find.Execute("apasiona debe ser mas cuidadosa", $true, $false, $false, $false, $false, $true, 1, $false, "apasiona debe ser más cuidadosa", 2)

# Word keeps the "_GoBack" bookmark pinned to the spot of the most recent
# edit; re-anchor it right after "más" (where the mas->más correction ended).
$r = $d.Content
$r.Find.Execute("más cuidadosa", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$goBackPos = $r.Start + 3
$goBackRange = $d.Range($goBackPos, $goBackPos)
$d.Bookmarks.Add("_GoBack", $goBackRange)
